# Update cryptocurrency price/volume data per commit
# (Mon Feb 26 07:00:12 UTC 2024 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.519.79"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.111.45"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "387.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0861"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "3.597.96"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "3.127.21"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").Value = "51.601.59"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  +6.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "0.0₃0969"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.166"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.99%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0479"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.21%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.68%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("E47").Value = "  +4.33%  "
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "2.080.96"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.925"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.14%  "
